# Remove the rows for MAGs GUT27127, GUT41097, GUT67224 and GUT80232.
# (Their predictions were superseded / the repo+output path was fixed, so
# those four genomes no longer belong in this sheet.) Deleting from the
# bottom up keeps the remaining row numbers stable while we work.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Rows.Item(11).EntireRow.Delete()   # even_MAG-GUT80232.fa
$ws.Rows.Item(7).EntireRow.Delete()    # even_MAG-GUT67224.fa
$ws.Rows.Item(4).EntireRow.Delete()    # even_MAG-GUT41097.fa
$ws.Rows.Item(3).EntireRow.Delete()    # even_MAG-GUT27127.fa
